$wb = $excel.ActiveWorkbook

# --- Sheet "stopwords" (sheet 1): append three new stopword rows ---
$wsStop = $wb.Worksheets.Item("stopwords")

$wsStop.Cells.Item(168, 1).Value = "use"
$wsStop.Cells.Item(168, 2).Value = "miscellaneous"

$wsStop.Cells.Item(169, 1).Value = "used"
$wsStop.Cells.Item(169, 2).Value = "miscellaneous"

$wsStop.Cells.Item(170, 1).Value = "using"
$wsStop.Cells.Item(170, 2).Value = "miscellaneous"

# --- View state: "stopwords" becomes the active sheet/tab ---
$wsStop.Activate()

$winStop = $excel.ActiveWindow
$winStop.ScrollRow = 154
$winStop.ScrollColumn = 1
$wsStop.Range("H173").Select()

# --- Sheet "keep in text" (sheet 2): update its (now inactive) selection ---
$wsKeep = $wb.Worksheets.Item("keep in text")
$wsKeep.Range("D7").Select()

# Re-activate "stopwords" so it is the tab shown when the workbook is reopened
$wsStop.Activate()
